# Swap the two theme colour palettes used by the deck:
#   theme1.xml ("Default" colours)            <->  theme2.xml ("Simple Light" / POSS colours)
#
# The presentation's slide master (and therefore every slide) is linked to
# ppt/theme/theme2.xml, which is the only theme colour scheme the
# PowerPoint object model exposes for editing (Master/NotesMaster/Slide all
# resolve to the same active ThemeColorScheme). We overwrite its twelve
# colours with the palette that used to live in theme1.xml ("Default"),
# completing the swap from the master's point of view.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$scheme = $m.Theme.ThemeColorScheme

# COM RGB values are packed as 0x00BBGGRR, so for a target hex "RRGGBB"
# the literal we assign is 0x00BBGGRR (bytes reversed).
$scheme.Item(1).RGB  = 0x00000000   # dk1      -> 000000
$scheme.Item(2).RGB  = 0x00FFFFFF   # lt1      -> FFFFFF
$scheme.Item(3).RGB  = 0x00588115   # dk2      -> 158158
$scheme.Item(4).RGB  = 0x00F3F3F3   # lt2      -> F3F3F3
$scheme.Item(5).RGB  = 0x00C78D05   # accent1  -> 058DC7
$scheme.Item(6).RGB  = 0x0032B450   # accent2  -> 50B432
$scheme.Item(7).RGB  = 0x001B56ED   # accent3  -> ED561B
$scheme.Item(8).RGB  = 0x0000EFED   # accent4  -> EDEF00
$scheme.Item(9).RGB  = 0x00E5CB24   # accent5  -> 24CBE5
$scheme.Item(10).RGB = 0x0072E564   # accent6  -> 64E572
$scheme.Item(11).RGB = 0x00CC0022   # hlink    -> 2200CC
$scheme.Item(12).RGB = 0x008B1A55   # folHlink -> 551A8B
